$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update 최종점수 (K) and MACRO_SCORE (N) values for rows 2-5
$ws.Range("K2").Value = 56.2
$ws.Range("N2").Value = 54.02451352198364

$ws.Range("K3").Value = 48.6
$ws.Range("N3").Value = 54.02451352198364

$ws.Range("K4").Value = 48.4
$ws.Range("N4").Value = 54.02451352198364

$ws.Range("K5").Value = 47.2
$ws.Range("N5").Value = 54.02451352198364
